$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the other header cells (row 1, e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record values for every data row (2 through 42)
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 100  # AD = Wins
    $ws.Cells.Item($row, 31).Value = 62   # AE = Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF = Ties
}
